$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; this shifts current rows 15-24 down to 16-25
# and inherits formatting (incl. the yyyy-mm-dd date format on column D) from
# the surrounding rows.
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with a fresh weekly record (same market/product
# grouping as the row below it, with an updated date and volume).
$ws.Cells.Item(15, 1).Value = 4
$ws.Cells.Item(15, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(15, 3).Value = "Los Lagos"
$ws.Cells.Item(15, 4).Value = 44894
$ws.Cells.Item(15, 5).Value = 10
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100107
$ws.Cells.Item(15, 8).Value = "Otros"
$ws.Cells.Item(15, 9).Value = 100107002
$ws.Cells.Item(15, 10).Value = "Chirimoya"
$ws.Cells.Item(15, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 200
$ws.Cells.Item(15, 14).Value = 22000
$ws.Cells.Item(15, 15).Value = 22500
$ws.Cells.Item(15, 16).Value = 22250
$ws.Cells.Item(15, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(15, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 19).Value = 2781
$ws.Cells.Item(15, 20).Value = 8
